# Applies the scheduled-runner price/profit refresh to the leve-flipping sheets.
# For each touched row, H/I/J/K/L are the live market-price columns and M/N are
# the derived NQ/HQ profit columns; values below mirror the refreshed dataset.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7422.923
$ws.Range("I40").Value = 2832
$ws.Range("K40").Value = 2832
$ws.Range("M40").Value = -2657

$ws.Range("H69").Value = 5800
$ws.Range("I69").Value = 5500
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 16500
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -15626
$ws.Range("N69").Value = -19748

$ws.Range("H72").Value = 5800
$ws.Range("I72").Value = 5500
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 49500
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -45132
$ws.Range("N72").Value = -62736

$ws.Range("H86").Value = 9999
$ws.Range("I86").Value = 9999
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 9999
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -8876
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 9999
$ws.Range("I89").Value = 9999
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 49995
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -44379
$ws.Range("N89").ClearContents()

$ws.Range("H108").Value = 37995
$ws.Range("J108").Value = 37995
$ws.Range("L108").Value = 37995
$ws.Range("N108").Value = -45675

$ws.Range("H137").Value = 3010.2666
$ws.Range("I137").Value = 1627
$ws.Range("K137").Value = 4881
$ws.Range("M137").Value = -2331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 77.6
$ws.Range("I5").Value = 84.5
$ws.Range("J5").Value = 73
$ws.Range("K5").Value = 84.5
$ws.Range("L5").Value = 73
$ws.Range("M5").Value = 27.5
$ws.Range("N5").Value = -297

$ws.Range("H32").Value = 9126.531
$ws.Range("I32").Value = 8068.4
$ws.Range("K32").Value = 8068.4
$ws.Range("M32").Value = -7781.4

$ws.Range("H61").Value = 3444
$ws.Range("I61").Value = 1925.6666
$ws.Range("K61").Value = 1925.6666
$ws.Range("M61").Value = -1713.6666

$ws.Range("H132").Value = 1813.4286
$ws.Range("I132").Value = 1813.4286
$ws.Range("K132").Value = 5440.2858
$ws.Range("M132").Value = -2910.2858

$ws.Range("H136").Value = 3444
$ws.Range("I136").Value = 1925.6666
$ws.Range("K136").Value = 5776.9998
$ws.Range("M136").Value = -3226.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 77.6
$ws.Range("I4").Value = 84.5
$ws.Range("J4").Value = 73
$ws.Range("K4").Value = 84.5
$ws.Range("L4").Value = 73
$ws.Range("M4").Value = 30.5
$ws.Range("N4").Value = -303

$ws.Range("H22").Value = 587.5
$ws.Range("I22").Value = 683.5
$ws.Range("J22").Value = 299.5
$ws.Range("K22").Value = 683.5
$ws.Range("L22").Value = 299.5
$ws.Range("M22").Value = -510.5
$ws.Range("N22").Value = -645.5

$ws.Range("H36").Value = 1846.25
$ws.Range("I36").Value = 143
$ws.Range("J36").Value = 3549.5
$ws.Range("K36").Value = 143
$ws.Range("L36").Value = 3549.5
$ws.Range("M36").Value = 391
$ws.Range("N36").Value = -4617.5

$ws.Range("H105").Value = 4873.6665
$ws.Range("J105").Value = 6499
$ws.Range("L105").Value = 6499
$ws.Range("N105").Value = -9993

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 838.75
$ws.Range("J22").Value = 690
$ws.Range("L22").Value = 690
$ws.Range("N22").Value = -1390

$ws.Range("H31").Value = 3050.0667
$ws.Range("I31").Value = 2865.6155
$ws.Range("K31").Value = 2865.6155
$ws.Range("M31").Value = -2570.6155

$ws.Range("H34").Value = 3050.0667
$ws.Range("I34").Value = 2865.6155
$ws.Range("K34").Value = 2865.6155
$ws.Range("M34").Value = -2663.6155

$ws.Range("H132").Value = 7638.3477
$ws.Range("I132").Value = 2975.5386
$ws.Range("J132").Value = 13700
$ws.Range("K132").Value = 8926.6158
$ws.Range("L132").Value = 41100
$ws.Range("M132").Value = -6396.6158
$ws.Range("N132").Value = -46160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 142.33333
$ws.Range("I38").Value = 66
$ws.Range("K38").Value = 198
$ws.Range("M38").Value = 149

$ws.Range("H75").Value = 867.8889
$ws.Range("I75").Value = 866.5
$ws.Range("K75").Value = 2599.5
$ws.Range("M75").Value = -1601.5

$ws.Range("H78").Value = 867.8889
$ws.Range("I78").Value = 866.5
$ws.Range("K78").Value = 7798.5
$ws.Range("M78").Value = -2806.5

$ws.Range("H114").Value = 398.85715
$ws.Range("J114").Value = 600
$ws.Range("L114").Value = 1800
$ws.Range("N114").Value = -8308

$ws.Range("H131").Value = 1382.0526
$ws.Range("I131").Value = 976.6
$ws.Range("J131").Value = 1832.5555
$ws.Range("K131").Value = 2929.8
$ws.Range("L131").Value = 5497.666499999999
$ws.Range("M131").Value = 2110.2
$ws.Range("N131").Value = -15577.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33336002
$ws.Range("I70").Value = 33336002
$ws.Range("K70").Value = 33336002
$ws.Range("M70").Value = -33335732

$ws.Range("H73").Value = 33336002
$ws.Range("I73").Value = 33336002
$ws.Range("K73").Value = 33336002
$ws.Range("M73").Value = -33335066

$ws.Range("H102").Value = 1070.6666
$ws.Range("I102").Value = 1070.6666
$ws.Range("K102").Value = 1070.6666
$ws.Range("M102").Value = 551.3334

$ws.Range("H126").Value = 280866800
$ws.Range("I126").Value = 280866800
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 842600400
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -842597930
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1167
$ws.Range("J16").Value = 980
$ws.Range("L16").Value = 980
$ws.Range("N16").Value = -1320

$ws.Range("H22").Value = 2380.3572
$ws.Range("I22").Value = 745
$ws.Range("J22").Value = 2826.3635
$ws.Range("K22").Value = 745
$ws.Range("L22").Value = 2826.3635
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -3416.3635

$ws.Range("H27").Value = 2380.3572
$ws.Range("I27").Value = 745
$ws.Range("J27").Value = 2826.3635
$ws.Range("K27").Value = 745
$ws.Range("L27").Value = 2826.3635
$ws.Range("M27").Value = -638
$ws.Range("N27").Value = -3040.3635

$ws.Range("H46").Value = 1428.4286
$ws.Range("I46").Value = 1399.5
$ws.Range("J46").Value = 1440
$ws.Range("K46").Value = 1399.5
$ws.Range("L46").Value = 1440
$ws.Range("M46").Value = -1211.5
$ws.Range("N46").Value = -1816

$ws.Range("H132").Value = 5999.2
$ws.Range("J132").Value = 10999.5
$ws.Range("L132").Value = 32998.5
$ws.Range("N132").Value = -38058.5

$ws.Range("H136").Value = 4392.9
$ws.Range("I136").Value = 3989.8572
$ws.Range("J136").Value = 5333.3335
$ws.Range("K136").Value = 11969.5716
$ws.Range("L136").Value = 16000.0005
$ws.Range("M136").Value = -9419.5716
$ws.Range("N136").Value = -21100.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9490

$ws.Range("H132").Value = 3850.3845
$ws.Range("I132").Value = 3850.3845
$ws.Range("K132").Value = 11551.1535
$ws.Range("M132").Value = -9021.1535
